$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 314
$ws1.Range("F5").Value = 207
$ws1.Range("F6").Value = 356
$ws1.Range("F8").Value = 2204
$ws1.Range("F9").Value = 375
$ws1.Range("F10").Value = 5396
$ws1.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg"
$ws1.Range("F12").Value = 358

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 314
$ws4.Range("F6").Value = 207
$ws4.Range("F7").Value = 356
$ws4.Range("F11").Value = 2204
$ws4.Range("F12").Value = 375
$ws4.Range("F13").Value = 5396
$ws4.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg"
$ws4.Range("F15").Value = 358
